$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- 1. Queue lock: flip "Save" flag (column E) off for rows 44-67 ---
for ($r = 44; $r -le 67; $r++) {
    $ws.Range("E$r").Value = $False
}

# --- 2. Re-format rows 76 & 77 so they match the plain (unshaded) rows ---
# Row 78 already carries the "no special formatting" look we need to copy
# onto most of rows 76/77; A76 keeps the lighter "s=1" look used elsewhere
# (e.g. A2) instead of the red-highlighted "s=3" look.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A76").PasteSpecial(-4122) | Out-Null

$ws.Range("A78").Copy() | Out-Null
$ws.Range("A77").PasteSpecial(-4122) | Out-Null

$ws.Range("B78").Copy() | Out-Null
$ws.Range("B76:B77").PasteSpecial(-4122) | Out-Null

$ws.Range("G78").Copy() | Out-Null
$ws.Range("G76:G77").PasteSpecial(-4122) | Out-Null

$ws.Range("H78").Copy() | Out-Null
$ws.Range("H76:H77").PasteSpecial(-4122) | Out-Null

$ws.Range("I78").Copy() | Out-Null
$ws.Range("I76:I77").PasteSpecial(-4122) | Out-Null

$ws.Range("J78").Copy() | Out-Null
$ws.Range("J76:J77").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 3. Move/restore the saved cursor position on this sheet ---
$ws.Range("H78").Select() | Out-Null
